# Removed "## " comment-prefixes from the R-console-output text blocks in
# anom-talk.pptx (slides 9, 18, 19, 21, 23, 29, 32).
#
# The "## " markers are literal text baked into the <a:t> runs (it is not a
# bullet or list marker), so each affected paragraph's text is rewritten in
# place. We address sub-ranges of each shape's TextRange with
# .Characters(start, length) (1-based, where multi-paragraph shapes burn one
# extra character per paragraph break) so that runs/paragraph breaks that do
# not change (e.g. the literal newlines already embedded inside a single
# run's text, or a trailing unrelated paragraph) are left completely alone.
# Where a shape has more than one paragraph being edited, the later
# paragraph is written first so the not-yet-updated earlier offsets stay
# valid.

$nl = [char]10

$p = $ppt.ActivePresentation

# --- Slide 9: "## png " / "##   2" -> "png " / "  2 " (2nd paragraph only) ---
$s = $p.Slides.Item(9)
$shp = $s.Shapes.Item(1)
$tr = $shp.TextFrame.TextRange
$tr.Characters(2, 14).Text = "png " + $nl + "  2 "

# --- Slide 18: "## Position 1" / data block ---
$s = $p.Slides.Item(18)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange
$tr.Characters(15, 103).Text = "-0.02386 -0.02853" + $nl + "-0.03001 -0.00428" + $nl + "-0.03623 -0.04222" + $nl + "-0.00144 -0.06466" + $nl + "0.00944 -0.00163"
$tr.Characters(1, 13).Text = "Position 1"

# --- Slide 19: "## Position 2" / data block ---
$s = $p.Slides.Item(19)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange
$tr.Characters(15, 97).Text = "-0.02014 -0.02725" + $nl + "0.02268 -0.03323" + $nl + "0.03661 0.04378" + $nl + "0.05562 0.00977" + $nl + "0.05641 0.01816"
$tr.Characters(1, 13).Text = "Position 2"

# --- Slide 21: summary statistics table ---
$s = $p.Slides.Item(21)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange
$tr.Characters(1, 188).Text = "  position    avg stdev" + $nl + "1        1 -0.022 0.023" + $nl + "2        2  0.016 0.033" + $nl + "3        3  0.006 0.029" + $nl + "4        4  0.065 0.021" + $nl + "5        5  0.008 0.026" + $nl + "6        6 -0.013 0.016"

# --- Slide 23: overall mean / pooled sd table ---
$s = $p.Slides.Item(23)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange
$tr.Characters(1, 31).Text = "   avg    sp" + $nl + "1 0.01 0.025"

# --- Slide 29: c_sections births table ---
$s = $p.Slides.Item(29)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange
$tr.Characters(1, 202).Text = "  group c_sections births" + $nl + "1    1A        150    923" + $nl + "2    1K         45    298" + $nl + "3    1B         34    170" + $nl + "4    1D         18    132" + $nl + "5    3I         20    106" + $nl + "6    3M         12    105"

# --- Slide 32: clinic cat_scans members table (leave the trailing prose paragraph alone) ---
$s = $p.Slides.Item(32)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange
$tr.Characters(1, 179).Text = "  clinic cat_scans members" + $nl + "1      1        50  26.838" + $nl + "2      2        71  26.895" + $nl + "3      3        41  26.142" + $nl + "4      4        62  25.907" + $nl + "5      5        89  26.565"
